$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1001423.8
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 1251529.8
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 3754589.4
$ws.Range("M17").Value = -2832
$ws.Range("N17").Value = -3754925.4

# Row 28
$ws.Range("H28").Value = 1359
$ws.Range("I28").Value = 1055.0834
$ws.Range("K28").Value = 1055.0834
$ws.Range("M28").Value = -570.0834

# Row 64
$ws.Range("H64").Value = 4312.4375
$ws.Range("I64").Value = 3461.5386
$ws.Range("K64").Value = 3461.5386
$ws.Range("M64").Value = -3213.5386

# Row 67
$ws.Range("H67").Value = 4312.4375
$ws.Range("I67").Value = 3461.5386
$ws.Range("K67").Value = 3461.5386
$ws.Range("M67").Value = -2603.5386

# Row 111
$ws.Range("H111").Value = 1125.3334
$ws.Range("J111").Value = 1306.2
$ws.Range("L111").Value = 3918.6
$ws.Range("N111").Value = -10052.6

# Row 127
$ws.Range("H127").Value = 673
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()


$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2215.3152
$ws.Range("I32").Value = 2198.6807
$ws.Range("K32").Value = 2198.6807
$ws.Range("M32").Value = -1911.6807

# Row 40
$ws.Range("H40").Value = 28052.334
$ws.Range("I40").Value = 15985
$ws.Range("J40").Value = 30465.8
$ws.Range("K40").Value = 15985
$ws.Range("L40").Value = 30465.8
$ws.Range("M40").Value = -15809
$ws.Range("N40").Value = -30817.8

# Row 45
$ws.Range("H45").Value = 6001.0967
$ws.Range("I45").Value = 7849
$ws.Range("J45").Value = 3442.4614
$ws.Range("K45").Value = 7849
$ws.Range("L45").Value = 3442.4614
$ws.Range("M45").Value = -7472
$ws.Range("N45").Value = -4196.4614

# Row 124
$ws.Range("H124").Value = 49585.8
$ws.Range("J124").Value = 49585.8
$ws.Range("L124").Value = 49585.8
$ws.Range("N124").Value = -59405.8

# Row 132
$ws.Range("H132").Value = 2459.4285
$ws.Range("I132").Value = 2459.4285
$ws.Range("K132").Value = 7378.2855
$ws.Range("M132").Value = -4848.2855


$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 1362.8
$ws.Range("I99").Value = 1362.8
$ws.Range("K99").Value = 1362.8
$ws.Range("M99").Value = 135.2


$ws = $wb.Worksheets.Item("CRP")
# Row 8
$ws.Range("H8").Value = 839
$ws.Range("J8").Value = 839
$ws.Range("L8").Value = 839
$ws.Range("N8").Value = -1119

# Row 107
$ws.Range("H107").Value = 2644.2173
$ws.Range("I107").Value = 2434
$ws.Range("K107").Value = 2434
$ws.Range("M107").Value = -514


$ws = $wb.Worksheets.Item("CUL")
# Row 45
$ws.Range("H45").Value = 12500
$ws.Range("J45").Value = 12500
$ws.Range("L45").Value = 37500
$ws.Range("N45").Value = -38564

# Row 70
$ws.Range("H70").Value = 4790
$ws.Range("I70").Value = 1661.6666
$ws.Range("K70").Value = 4984.9998
$ws.Range("M70").Value = -4669.9998

# Row 73
$ws.Range("H73").Value = 4790
$ws.Range("I73").Value = 1661.6666
$ws.Range("K73").Value = 4984.9998
$ws.Range("M73").Value = -3892.9998

# Row 121
$ws.Range("H121").Value = 20835252
$ws.Range("I121").Value = 83333620
$ws.Range("K121").Value = 250000860
$ws.Range("M121").Value = -249999550

# Row 132
$ws.Range("H132").Value = 2441.1365
$ws.Range("J132").Value = 2747
$ws.Range("L132").Value = 24723
$ws.Range("N132").Value = -29783


$ws = $wb.Worksheets.Item("GSM")
# Row 19
$ws.Range("H19").Value = 360
$ws.Range("I19").Value = 360
$ws.Range("K19").Value = 360
$ws.Range("M19").Value = -72

# Row 33
$ws.Range("H33").Value = 8500
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 8500
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 8500
$ws.Range("M33").ClearContents()
$ws.Range("N33").Value = -9004

# Row 43
$ws.Range("H43").Value = 12773.714
$ws.Range("I43").Value = 6569.5
$ws.Range("K43").Value = 6569.5
$ws.Range("M43").Value = -6418.5

# Row 96
$ws.Range("H96").Value = 29761
$ws.Range("J96").Value = 29761
$ws.Range("L96").Value = 29761
$ws.Range("N96").Value = -35253

# Row 102
$ws.Range("H102").Value = 2776.7307
$ws.Range("I102").Value = 2742
$ws.Range("J102").Value = 2922.6
$ws.Range("K102").Value = 2742
$ws.Range("L102").Value = 2922.6
$ws.Range("M102").Value = -1120
$ws.Range("N102").Value = -6166.6

# Row 113
$ws.Range("H113").Value = 4820
$ws.Range("I113").Value = 6498.75
$ws.Range("K113").Value = 6498.75
$ws.Range("M113").Value = -4328.75

# Row 132
$ws.Range("H132").Value = 7983.8184
$ws.Range("I132").Value = 6628.7144
$ws.Range("J132").Value = 10355.25
$ws.Range("K132").Value = 19886.1432
$ws.Range("L132").Value = 31065.75
$ws.Range("M132").Value = -17356.1432
$ws.Range("N132").Value = -36125.75


$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 345.5
$ws.Range("I16").Value = 345.5
$ws.Range("K16").Value = 345.5
$ws.Range("M16").Value = -175.5

# Row 55
$ws.Range("H55").Value = 258.84616
$ws.Range("I55").Value = 265.57144
$ws.Range("J55").Value = 251
$ws.Range("K55").Value = 265.57144
$ws.Range("L55").Value = 251
$ws.Range("M55").Value = -92.57144
$ws.Range("N55").Value = -597

# Row 61
$ws.Range("H61").Value = 45347.566
$ws.Range("I61").Value = 51964.85
$ws.Range("K61").Value = 51964.85
$ws.Range("M61").Value = -51762.85

# Row 68
$ws.Range("H68").Value = 2088
$ws.Range("I68").Value = 1117.3334
$ws.Range("K68").Value = 1117.3334
$ws.Range("M68").Value = -368.3334

# Row 71
$ws.Range("H71").Value = 2088
$ws.Range("I71").Value = 1117.3334
$ws.Range("K71").Value = 5586.666999999999
$ws.Range("M71").Value = -1842.666999999999

# Row 94
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

# Row 100
$ws.Range("H100").Value = 7000
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 7000
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 7000
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -8082

# Row 113
$ws.Range("H113").Value = 45347.566
$ws.Range("I113").Value = 51964.85
$ws.Range("K113").Value = 51964.85
$ws.Range("M113").Value = -49794.85

# Row 122
$ws.Range("H122").Value = 158036.08
$ws.Range("I122").Value = 203546.9
$ws.Range("K122").Value = 610640.7
$ws.Range("M122").Value = -608190.7

# Row 140
$ws.Range("H140").Value = 94427
$ws.Range("J140").Value = 94427
$ws.Range("L140").Value = 94427
$ws.Range("N140").Value = -104787


$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 48000
$ws.Range("I2").Value = 47000
$ws.Range("J2").Value = 50000
$ws.Range("K2").Value = 47000
$ws.Range("L2").Value = 50000
$ws.Range("M2").Value = -46888
$ws.Range("N2").Value = -50224

# Row 14
$ws.Range("H14").Value = 1206996
$ws.Range("I14").Value = 3001000
$ws.Range("J14").Value = 10993.333
$ws.Range("K14").Value = 3001000
$ws.Range("L14").Value = 10993.333
$ws.Range("M14").Value = -3000832
$ws.Range("N14").Value = -11329.333

# Row 69
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

# Row 72
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

# Row 107
$ws.Range("H107").Value = 14286902
$ws.Range("I107").Value = 1156.2916
$ws.Range("K107").Value = 3468.8748
$ws.Range("M107").Value = -1548.8748

# Row 138
$ws.Range("H138").Value = 94998
$ws.Range("J138").Value = 94998
$ws.Range("L138").Value = 94998
$ws.Range("N138").Value = -105278

